# typeset_drafts/192200_radio_for_all.docx
#
# The source Markdown had three inline images whose link syntax leaked
# into the prose as literal runs of text, e.g. `...](images/foo.jpg)`
# (with the caption / alt text in the preceding run(s)). Replace each of
# those remnants with a plain "[INSERT FIGURE n.n NEAR HERE]" placeholder
# so the typeset draft reads cleanly and an editor can drop the real
# figure in later.

$d = $word.ActiveDocument

function Replace-ParagraphText($marker, $newText) {
    # Find the paragraph whose text contains $marker, then overwrite the
    # paragraph's whole content (which may be spread across several runs)
    # with a single run holding $newText, leaving the paragraph mark (and
    # anything that follows it, e.g. bookmarks) untouched.
    foreach ($p in $d.Paragraphs) {
        $r = $p.Range
        if ($r.Text.Contains($marker)) {
            $r.InsertBefore($newText)
            $old = $d.Range($r.Start + $newText.Length, $r.End)
            $old.Delete()
            return $true
        }
    }
    return $false
}

# --- Figure 39.1 ---------------------------------------------------
# "<curly-quote>It is a mistake ... as yet.<curly-quote>](images/radio_for_all.jpg)"
Replace-ParagraphText "(images/radio_for_all.jpg)" "[INSERT FIGURE 39.1 NEAR HERE]" | Out-Null

# --- Figure 39.2 ---------------------------------------------------
# "](images/radio_for_all2.jpg)"
Replace-ParagraphText "(images/radio_for_all2.jpg)" "[INSERT FIGURE 39.2 NEAR HERE]" | Out-Null

# --- Figure 39.3 ---------------------------------------------------
# "The power plant at the radio broadcasting station ... away.](images/radio_for_all3.jpg)"
Replace-ParagraphText "(images/radio_for_all3.jpg)" "[INSERT FIGURE 39.3 NEAR HERE]" | Out-Null

Write-Output "figure placeholders inserted"
